$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.986.88'
$ws.Range("E2").Value = '  +4.37%  '
$ws.Range("D3").Value = '3.252.79'
$ws.Range("E3").Value = '  +2.21%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '579.26'
$ws.Range("E5").Value = '  +3.05%  '
$ws.Range("D6").Value = '177.53'
$ws.Range("E6").Value = '  +3.47%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.606'
$ws.Range("E7").Value = '  +0.84%  '
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '3.249.45'
$ws.Range("E9").Value = '  +2.18%  '
$ws.Range("E10").Value = '  +4.60%  '
$ws.Range("D11").Value = '6.72'
$ws.Range("E11").Value = '  +1.33%  '
$ws.Range("D12").Value = '0.409'
$ws.Range("E12").Value = '  +3.19%  '
$ws.Range("D13").Value = '3.811.88'
$ws.Range("E13").Value = '  +2.11%  '
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("D15").Value = '28.16'
$ws.Range("E15").Value = '  +1.82%  '
$ws.Range("D16").Value = '66.958.64'
$ws.Range("E16").Value = '  +4.34%  '
$ws.Range("D17").Value = '0.0000168'
$ws.Range("E17").Value = '  +3.31%  '
$ws.Range("D18").Value = '3.249.11'
$ws.Range("E18").Value = '  +2.22%  '
$ws.Range("D19").Value = '5.81'
$ws.Range("E19").Value = '  +2.14%  '
$ws.Range("D20").Value = '13.40'
$ws.Range("E20").Value = '  +2.10%  '
$ws.Range("D21").Value = '369.75'
$ws.Range("E21").Value = '  +4.84%  '
$ws.Range("D22").Value = '7.60'
$ws.Range("E22").Value = '  +5.82%  '
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").Value = '70.80'
$ws.Range("E24").Value = '  +1.86%  '
$ws.Range("D25").Value = '0.509'
$ws.Range("E25").Value = '  +0.97%  '
$ws.Range("D26").Value = '3.388.00'
$ws.Range("E26").Value = '  +2.10%  '
$ws.Range("D27").Value = '0.0000119'
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("D28").Value = '9.79'
$ws.Range("E28").Value = '  +3.14%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("E31").Value = '  +4.71%  '
$ws.Range("D32").Value = '5.65'
$ws.Range("E32").Value = '  -0.94%  '
$ws.Range("D33").Value = '22.56'
$ws.Range("E33").Value = '  +1.73%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("D35").Value = '174.46'
$ws.Range("E35").Value = '  +11.19%  '
$ws.Range("E36").Value = '  +2.77%  '
$ws.Range("D37").Value = '6.77'
$ws.Range("E37").Value = '  +1.98%  '
$ws.Range("D38").Value = '1.52'
$ws.Range("E38").Value = '  +5.34%  '
$ws.Range("D39").Value = '0.857'
$ws.Range("E39").Value = '  +6.45%  '
$ws.Range("E40").Value = '  +9.68%  '
$ws.Range("D41").Value = '26.86'
$ws.Range("E41").Value = '  +3.40%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = '2.57'
$ws.Range("E42").Value = '  +2.75%  '
$ws.Range("D43").Value = '2.753.92'
$ws.Range("E43").Value = '  +2.98%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '6.41'
$ws.Range("E44").Value = '  +6.16%  '
$ws.Range("D45").Value = '4.31'
$ws.Range("E45").Value = '  +3.59%  '
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").Value = '340.13'
$ws.Range("E46").Value = '  +3.83%  '
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").Value = '40.31'
$ws.Range("E47").Value = '  +4.50%  '
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").Value = '0.0674'
$ws.Range("E48").Value = '  +3.32%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '24.74'
$ws.Range("E49").Value = '  +4.48%  '
$ws.Range("E51").Value = '  +2.41%  '
